# "Updated CVDs for the month"
# Zero out the Jul-Dec / Q3 / Q4 / FY "Commit/Forecast" CVD figures for the
# Milwaukee Pmc Hq Wisconsin "Professional Voluntary Turnover" row, and clear
# the now-stale Jul "Internal Fill Rate" / "Professional Voluntary Turnover"
# placeholder zeros (O5 on Milwaukee, O4 on South Beloit) back to blank cells.

$wb = $excel.ActiveWorkbook

$wsMilwaukee = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$wsMilwaukee.Range("O4:W4").Value = 0
$wsMilwaukee.Range("O5").ClearContents()

$wsSouthBeloit = $wb.Worksheets.Item("South Beloit Gardner St Illino")
$wsSouthBeloit.Range("O4").ClearContents()
